# Update the "Date insertion" (column H) values for the first 7 films.
# The values are date-like strings that must stay as literal text, so we
# force a Text number format before assigning the value (otherwise Excel's
# automatic type conversion would turn them into date serial numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "H2" = "06/01/2019"
    "H3" = "17/06/2018"
    "H4" = "11/05/2019"
    "H5" = "25/08/2018"
    "H6" = "23/05/2018"
    "H7" = "28/01/2018"
    "H8" = "20/01/2019"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
